$wb = $excel.ActiveWorkbook

# --- 1. Update status text everywhere it appears ("Ready for handoff" -> "In Translation") ---

# Overview sheet: zh-cn (E) / de-de (F) status cells for both data rows
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for both data rows
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "In Translation"
$wsZh.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for both data rows
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "In Translation"
$wsDe.Range("C3").Value = "In Translation"

# --- 2. Resize the status columns that used to be sized for "Ready for handoff" ---
# Overview sheet: columns E (zh-cn) and F (de-de)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn sheet: column C (Status)
$wsZh.Columns.Item(3).ColumnWidth = 12.5

# de-de sheet: column C (Status)
$wsDe.Columns.Item(3).ColumnWidth = 12.5
